$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2:G2")
$rng.NumberFormat = "@"
$ws.Range("D2").Value = "261.55"
$ws.Range("E2").Value = "0.96%"
$ws.Range("G2").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("E3:G3")
$rng.NumberFormat = "@"
$ws.Range("E3").Value = "0.61%"
$ws.Range("G3").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D4:G4")
$rng.NumberFormat = "@"
$ws.Range("D4").Value = "4.710"
$ws.Range("E4").Value = "0.66%"
$ws.Range("G4").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D5:G5")
$rng.NumberFormat = "@"
$ws.Range("D5").Value = "0.06208"
$ws.Range("E5").Value = "2.88%"
$ws.Range("G5").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D6:G6")
$rng.NumberFormat = "@"
$ws.Range("D6").Value = "6.724"
$ws.Range("E6").Value = "0.74%"
$ws.Range("G6").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D7:G7")
$rng.NumberFormat = "@"
$ws.Range("D7").Value = "0.8499"
$ws.Range("E7").Value = "-1.17%"
$ws.Range("G7").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D8:G8")
$rng.NumberFormat = "@"
$ws.Range("D8").Value = "0.9105"
$ws.Range("E8").Value = "-1.09%"
$ws.Range("G8").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("E9:G9")
$rng.NumberFormat = "@"
$ws.Range("E9").Value = "0.82%"
$ws.Range("G9").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D10:G10")
$rng.NumberFormat = "@"
$ws.Range("D10").Value = "0.04680"
$ws.Range("E10").Value = "-8.03%"
$ws.Range("G10").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D11:G11")
$rng.NumberFormat = "@"
$ws.Range("D11").Value = "0.07096"
$ws.Range("E11").Value = "0.23%"
$ws.Range("G11").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D12:G12")
$rng.NumberFormat = "@"
$ws.Range("D12").Value = "0.03175"
$ws.Range("E12").Value = "3.25%"
$ws.Range("G12").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D13:G13")
$rng.NumberFormat = "@"
$ws.Range("D13").Value = "0.09065"
$ws.Range("E13").Value = "-0.66%"
$ws.Range("G13").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D14:G14")
$rng.NumberFormat = "@"
$ws.Range("D14").Value = "0.001534"
$ws.Range("E14").Value = "0.08%"
$ws.Range("G14").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D15:G15")
$rng.NumberFormat = "@"
$ws.Range("D15").Value = "0.0006142"
$ws.Range("E15").Value = "1.66%"
$ws.Range("G15").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D16:G16")
$rng.NumberFormat = "@"
$ws.Range("D16").Value = "0.006145"
$ws.Range("E16").Value = "0.07%"
$ws.Range("G16").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D17:G17")
$rng.NumberFormat = "@"
$ws.Range("D17").Value = "3.466"
$ws.Range("G17").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D18:G18")
$rng.NumberFormat = "@"
$ws.Range("D18").Value = "3.171"
$ws.Range("E18").Value = "-0.11%"
$ws.Range("G18").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D19:G19")
$rng.NumberFormat = "@"
$ws.Range("D19").Value = "2.178"
$ws.Range("E19").Value = "0.57%"
$ws.Range("G19").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("E20:G20")
$rng.NumberFormat = "@"
$ws.Range("E20").Value = "-0.55%"
$ws.Range("G20").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D21:G21")
$rng.NumberFormat = "@"
$ws.Range("D21").Value = "0.1300"
$ws.Range("E21").Value = "0.18%"
$ws.Range("G21").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D22:G22")
$rng.NumberFormat = "@"
$ws.Range("D22").Value = "4.118"
$ws.Range("E22").Value = "-0.43%"
$ws.Range("G22").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D23:G23")
$rng.NumberFormat = "@"
$ws.Range("D23").Value = "0.04232"
$ws.Range("E23").Value = "-0.04%"
$ws.Range("G23").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D24:G24")
$rng.NumberFormat = "@"
$ws.Range("D24").Value = "0.001216"
$ws.Range("E24").Value = "-0.13%"
$ws.Range("G24").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D25:G25")
$rng.NumberFormat = "@"
$ws.Range("D25").Value = "0.004134"
$ws.Range("E25").Value = "2.58%"
$ws.Range("G25").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("E26:G26")
$rng.NumberFormat = "@"
$ws.Range("E26").Value = "0.10%"
$ws.Range("G26").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D27:G27")
$rng.NumberFormat = "@"
$ws.Range("D27").Value = "0.0001617"
$ws.Range("E27").Value = "6.14%"
$ws.Range("G27").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("G28")
$rng.NumberFormat = "@"
$ws.Range("G28").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("G29")
$rng.NumberFormat = "@"
$ws.Range("G29").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("G30")
$rng.NumberFormat = "@"
$ws.Range("G30").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("G31")
$rng.NumberFormat = "@"
$ws.Range("G31").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("G32")
$rng.NumberFormat = "@"
$ws.Range("G32").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("G33")
$rng.NumberFormat = "@"
$ws.Range("G33").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("G34")
$rng.NumberFormat = "@"
$ws.Range("G34").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("G35")
$rng.NumberFormat = "@"
$ws.Range("G35").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("G36")
$rng.NumberFormat = "@"
$ws.Range("G36").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("G37")
$rng.NumberFormat = "@"
$ws.Range("G37").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("G38")
$rng.NumberFormat = "@"
$ws.Range("G38").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("G39")
$rng.NumberFormat = "@"
$ws.Range("G39").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D40:G40")
$rng.NumberFormat = "@"
$ws.Range("D40").Value = "0.03912"
$ws.Range("E40").Value = "1.71%"
$ws.Range("G40").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("E41:G41")
$rng.NumberFormat = "@"
$ws.Range("E41").Value = "0.02%"
$ws.Range("G41").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("E42:G42")
$rng.NumberFormat = "@"
$ws.Range("E42").Value = "2.66%"
$ws.Range("G42").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D43:G43")
$rng.NumberFormat = "@"
$ws.Range("D43").Value = "0.002184"
$ws.Range("E43").Value = "-0.73%"
$ws.Range("G43").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D44:G44")
$rng.NumberFormat = "@"
$ws.Range("D44").Value = "0.01352"
$ws.Range("E44").Value = "-11.55%"
$ws.Range("G44").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D45:G45")
$rng.NumberFormat = "@"
$ws.Range("D45").Value = "0.00005175"
$ws.Range("E45").Value = "1.56%"
$ws.Range("G45").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D46:G46")
$rng.NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("E46").Value = "0.09%"
$ws.Range("G46").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D47:G47")
$rng.NumberFormat = "@"
$ws.Range("D47").Value = "0.03592"
$ws.Range("E47").Value = "-34.14%"
$ws.Range("G47").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("D48:G48")
$rng.NumberFormat = "@"
$ws.Range("D48").Value = "0.1675"
$ws.Range("E48").Value = "26.79%"
$ws.Range("G48").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("E49:G49")
$rng.NumberFormat = "@"
$ws.Range("E49").Value = "0.09%"
$ws.Range("G49").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("E50:G50")
$rng.NumberFormat = "@"
$ws.Range("E50").Value = "0.09%"
$ws.Range("G50").Value = "22"
$rng.ClearFormats()

$rng = $ws.Range("G51")
$rng.NumberFormat = "@"
$ws.Range("G51").Value = "22"
$rng.ClearFormats()
